# Rotate the data rows 6-12 of the sheet one step "up" the cycle:
#   row7 -> row6, row8 -> row7, row9 -> row8, row10 -> row9,
#   row12 -> row10, row11 -> row12, row6 -> row11
# (i.e. a single 7-cycle permutation of whole rows; columns/layout untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns holding plain text in rows 6-12 (must stay text - force Text format so
# Excel's autodetection doesn't reinterpret date-looking strings like
# "2018-09-12" as real dates, or other strings as numbers).
$textCols = @("C","D","F","G","H","I","P","T","U","V","W","Y","Z","AA","AB","AT","AW","AX","AY")

# Columns holding real numbers.
$numCols = @("A","B","E","Q","R","S")

# Columns holding booleans.
$boolCols = @("AD","AE","AG")

$allCols = $numCols + $textCols + $boolCols
$rows = @(6,7,8,9,10,11,12)

# 1) Snapshot every cell in rows 6-12 BEFORE writing anything back, since the
#    write step below overwrites rows that other rows' new values are sourced from.
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $allCols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# 2) Destination row -> source row (content moves FROM source TO destination).
$mapping = @{
    6  = 7
    7  = 8
    8  = 9
    9  = 10
    10 = 12
    11 = 6
    12 = 11
}

# 3) Make sure the text columns keep a Text number format on all affected rows,
#    so re-writing values such as "2018-09-12" does not get auto-converted to a
#    date serial number.
foreach ($r in $rows) {
    foreach ($c in $textCols) {
        $ws.Range("$c$r").NumberFormat = "@"
    }
}

# 4) Write each destination row's cells from the captured snapshot of its source row.
foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $allCols) {
        $ws.Range("$c$destRow").Value = $srcData[$c]
    }
}
